$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Paragraph 11 ("Between 2000 to 2006 Ecuador had 4 presidents...")
#    replace with the expanded / reworded text.
# ---------------------------------------------------------------------
$old11 = "Between 2000 to 2006 Ecuador had 4 presidents, society was extremely affected and the presidents could not hold their positions due to the severity of the economic crisis. "
$new11 = "Between 2000 to 2006, Ecuador had four presidents. Consequently, society was damaged and the presidents could not hold their representative positions in government due to the severity of the economic crisis and the pressure from the people."
$d.Content.Find.Execute($old11, $true, $false, $false, $false, $false, $true, 1, $false, $new11, 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Paragraph 12 ("In 2007 Rafael Correa runs for office and launches
#    a political proposal sifting power...") replace with expanded text.
# ---------------------------------------------------------------------
$old12 = "In 2007 Rafael Correa runs for office and launches a political proposal sifting power from the high class and investing resources in low income communities. The movement lasted 10 years and was able to reduce Ecuador" + [char]8217 + "s poverty rate, build schools, hospitals, roads and over all revive the damaged economy. "
$new12 = "In 2007 Rafael Correa runs for office and launches a political proposal to reduce power from the elite families and companies that controlled the country to then invest resources in low income communities providing quality education, road infrastructure, diversifying local and foreign investment to improve the way of life of people. The movement lasted 10 years and reduced Ecuador" + [char]8217 + "s population poverty rate from 64.4% in 2000 to 22.9% at the end of 2017. "
$d.Content.Find.Execute($old12, $true, $false, $false, $false, $false, $true, 1, $false, $new12, 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Insert the new paragraphs that follow paragraph 12 (the "In 2007..."
#    paragraph) and precede the old "(empty)" / "Sources:" paragraphs:
#      - Correa received severe criticism...
#      - (empty paragraph)
#      - CHART (bold)
#      - As life continues in Ecuador...
#      - Above all, this message resonates...
# ---------------------------------------------------------------------
$correaPara = $d.Paragraphs.Item(12)

$correaPara.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(13)
$p.Range.Text = "Correa received severe criticism from the conservative right in Ecuador. The majority of the media is owned by affluent individuals who were affected through the new policies. He was labeled as a communist or a socialist due to the progressive reforms implemented. "

$d.Paragraphs.Item(13).Range.InsertParagraphAfter()
# paragraph 14 stays empty

$d.Paragraphs.Item(14).Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(15)
$p.Range.Text = "CHART"
$p.Range.Bold = 1

$d.Paragraphs.Item(15).Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(16)
$p.Range.Text = "As life continues in Ecuador, people have live through the changes in these past 18 years since the financial crisis. In the same way, Ecuador" + [char]8217 + "s GDP has increased significantly in this past decade allowing the economy to stabilize. However, this has also created a division between family members who are either against or for Correa" + [char]8217 + "s movement. "

$d.Paragraphs.Item(16).Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(17)
$p.Range.Text = "Above all, this message resonates with political environments around the world. The media and some politicians have purposely divided people and placed them against each other by empowering extremist in both sides to create chaos and cut dialogue between individuals. This is how societal crisis start."

# ---------------------------------------------------------------------
# 4) Move the "_GoBack" bookmark from the last paragraph (URL) to a new
#    empty paragraph right after the "Above all..." paragraph, and add
#    a handful of break-only paragraphs before "Sources:" picks back up.
# ---------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$d.Paragraphs.Item(17).Range.InsertParagraphAfter()
$gobackPara = $d.Paragraphs.Item(18)
$d.Bookmarks.Add("_GoBack", $gobackPara.Range)

$d.Paragraphs.Item(18).Range.InsertParagraphAfter()
# paragraph 19 stays empty (matches the old empty paragraph that used to
# sit right before "Sources:")

$d.Paragraphs.Item(19).Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(20)
$p.Range.InsertBreak(6) | Out-Null

$d.Paragraphs.Item(20).Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(21)
$p.Range.InsertBreak(6) | Out-Null
$p.Range.InsertAfter(" ")
$p.Range.InsertAfter(" ")

# ---------------------------------------------------------------------
# 5) Tidy up the repeated "Poverty headcount ratio..." / "DataBank..."
#    paragraphs: collapse the split proofErr runs into single runs (no
#    textual change, just simplifying markup) -- handled automatically
#    by the Find/Replace below which re-writes those paragraphs whole.
# ---------------------------------------------------------------------
$povOld = "Poverty headcount ratio at national poverty lines (% of population)"
$povNew = "Poverty headcount ratio at national poverty lines (% of population)"
# Re-run Find/Replace (matches whole text) so that the proofErr-split runs
# collapse into a single run, for both occurrences.
$rng = $d.Content
$rng.Find.Execute($povOld, $true, $false, $false, $false, $false, $true, 1, $false, $povNew, 2) | Out-Null
$rng = $d.Content
$rng.Find.Execute($povOld, $true, $false, $false, $false, $false, $true, 1, $false, $povNew, 2) | Out-Null

$dbOld = "DataBankMicrodataData Catalog"
$dbNew = "DataBankMicrodataData Catalog"
$rng = $d.Content
$rng.Find.Execute($dbOld, $true, $false, $false, $false, $false, $true, 1, $false, $dbNew, 2) | Out-Null

# ---------------------------------------------------------------------
# 6) Append the new closing content after the migrationpolicy.org URL
#    paragraph: an empty paragraph followed by the "People seeing this
#    phenomena..." paragraph.
# ---------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$lastPara.Range.InsertParagraphAfter()
$lastParaIndex = $lastParaIndex + 1
# new empty paragraph stays empty

$d.Paragraphs.Item($lastParaIndex).Range.InsertParagraphAfter()
$lastParaIndex = $lastParaIndex + 1
$p = $d.Paragraphs.Item($lastParaIndex)
$p.Range.Text = "People seeing this phenomena, decided to come back to the country and by 2015 we see a great decline of migration out of the country. "
